$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = "18:31:39"
$ws.Range("D9").Value = "18:31:39"
$ws.Range("D10").Value = "18:31:39"
$ws.Range("D11").Value = "18:31:39"
$ws.Range("D12").Value = "18:46:39"
$ws.Range("D13").Value = "18:46:39"
$ws.Range("D14").Value = "18:46:39"
$ws.Range("D15").Value = "18:46:39"
$ws.Range("D16").Value = "19:11:39"
$ws.Range("D17").Value = "19:11:39"
$ws.Range("D18").Value = "19:11:39"
$ws.Range("D19").Value = "19:11:39"
$ws.Range("D20").Value = "19:36:39"
$ws.Range("D21").Value = "19:36:39"
$ws.Range("D22").Value = "19:36:39"
$ws.Range("D23").Value = "19:36:39"
$ws.Range("D24").Value = "19:53:18"
$ws.Range("D25").Value = "19:53:18"
$ws.Range("D26").Value = "19:53:18"
$ws.Range("D27").Value = "19:53:18"
$ws.Range("D28").Value = "20:08:18"
$ws.Range("D29").Value = "20:08:18"
$ws.Range("D30").Value = "20:08:18"
$ws.Range("D31").Value = "20:08:18"
$ws.Range("D32").Value = "20:33:18"
$ws.Range("D33").Value = "20:33:18"
$ws.Range("D34").Value = "20:33:18"
$ws.Range("D35").Value = "20:33:18"
$ws.Range("D36").Value = "20:58:18"
$ws.Range("D37").Value = "20:58:18"
$ws.Range("D38").Value = "20:58:18"
$ws.Range("D39").Value = "20:58:18"
$ws.Range("D40").Value = "21:14:57"
$ws.Range("D41").Value = "21:14:57"
$ws.Range("D42").Value = "21:14:57"
$ws.Range("D43").Value = "21:14:57"
$ws.Range("D44").Value = "21:29:57"
$ws.Range("D45").Value = "21:29:57"
$ws.Range("D46").Value = "21:29:57"
$ws.Range("D47").Value = "21:29:57"
$ws.Range("D48").Value = "21:54:57"
$ws.Range("D49").Value = "21:54:57"
$ws.Range("D50").Value = "21:54:57"
$ws.Range("D51").Value = "21:54:57"
$ws.Range("D52").Value = "22:19:57"
$ws.Range("D53").Value = "22:19:57"
$ws.Range("D54").Value = "22:19:57"
$ws.Range("D55").Value = "22:19:57"
$ws.Range("D56").Value = "22:36:36"
$ws.Range("D57").Value = "22:36:36"
$ws.Range("D58").Value = "22:36:36"
$ws.Range("D59").Value = "22:36:36"
$ws.Range("D60").Value = "22:51:36"
$ws.Range("D61").Value = "22:51:36"
$ws.Range("D62").Value = "22:51:36"
$ws.Range("D63").Value = "22:51:36"
$ws.Range("D64").Value = "23:16:36"
$ws.Range("D65").Value = "23:16:36"
$ws.Range("D66").Value = "23:16:36"
$ws.Range("D67").Value = "23:16:36"
$ws.Range("D68").Value = "23:41:36"
$ws.Range("D69").Value = "23:41:36"
$ws.Range("D70").Value = "23:41:36"
$ws.Range("D71").Value = "23:41:36"
$ws.Range("D72").Value = "23:58:15"
$ws.Range("D73").Value = "23:58:15"
$ws.Range("D74").Value = "23:58:15"
$ws.Range("D75").Value = "23:58:15"
